# Auto-generated edit script applying cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.651.21'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.304.16'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.06'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.83'
$ws.Range("E6").Value = '  -6.45%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.581'
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.298.87'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("E10").Value = '  -4.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.574'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.32'
$ws.Range("E12").Value = '  -4.54%  '
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '668.64'
$ws.Range("E14").Value = '  +4.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.847.42'
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.34'
$ws.Range("E16").Value = '  -2.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.795.36'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.315.23'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.82'
$ws.Range("E21").Value = '  -3.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.885'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.11'
$ws.Range("E23").Value = '  -5.23%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.36'
$ws.Range("E24").Value = '  +4.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.29'
$ws.Range("E25").Value = '  -2.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.85'
$ws.Range("E26").Value = '  -4.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("E27").Value = '  -6.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("E28").Value = '  -5.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.00'
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.36'
$ws.Range("E30").Value = '  -3.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.04'
$ws.Range("E31").Value = '  +1.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '590.60'
$ws.Range("E32").Value = '  -2.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.91'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("E34").Value = '  -2.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.734.53'
$ws.Range("E35").Value = '  -6.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.37'
$ws.Range("E37").Value = '  -12.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.27'
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.131'
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("E40").Value = '  -8.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.18'
$ws.Range("E41").Value = '  -4.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.10'
$ws.Range("E42").Value = '  -4.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0664'
$ws.Range("E43").Value = '  -5.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.28'
$ws.Range("E44").Value = '  -3.39%  '
$ws.Range("E45").Value = '  -4.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0405'
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.34'
$ws.Range("E50").Value = '  -1.73%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '128.81'
$ws.Range("E51").Value = '  +0.24%  '
